# Auto-generated script to apply numeric corrections to the Ravana_Profits market-data sheets.
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) has per-leve market price/profit
# columns (H-N) that were refreshed from an updated data pull. This sets the corrected
# values cell-by-cell, and clears cells that no longer apply (e.g. no HQ price/profit).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2500
$ws.Range("I64").Value = 2500
$ws.Range("K64").Value = 2500
$ws.Range("M64").Value = -2252
# Row 67
$ws.Range("H67").Value = 2500
$ws.Range("I67").Value = 2500
$ws.Range("K67").Value = 2500
$ws.Range("M67").Value = -1642
# Row 74
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
# Row 77
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
# Row 96
$ws.Range("H96").Value = 20000776
$ws.Range("I96").Value = 971
$ws.Range("K96").Value = 2913
$ws.Range("M96").Value = -1540
# Row 106
$ws.Range("H106").Value = 7887.5
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -9262
# Row 107
$ws.Range("H107").Value = 478.36365
$ws.Range("I107").Value = 506.2
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 506.2
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 1413.8
$ws.Range("N107").Value = -4040
# Row 112
$ws.Range("H112").Value = 1873.2858
$ws.Range("J112").Value = 1955.8462
$ws.Range("L112").Value = 5867.5386
$ws.Range("N112").Value = -8083.5386
# Row 113
$ws.Range("H113").Value = 5001500
$ws.Range("I113").Value = 10000000
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 10000000
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -9996746
$ws.Range("N113").Value = -9508
# Row 116
$ws.Range("H116").Value = 9849.143
$ws.Range("I116").Value = 9648.333000000001
$ws.Range("K116").Value = 9648.333000000001
$ws.Range("M116").Value = -6206.333000000001
# Row 132
$ws.Range("H132").Value = 1116.3572
$ws.Range("I132").Value = 1116.3572
$ws.Range("K132").Value = 3349.0716
$ws.Range("M132").Value = -819.0715999999998
# Row 137
$ws.Range("H137").Value = 2292.75
$ws.Range("I137").Value = 1630.875
$ws.Range("J137").Value = 2954.625
$ws.Range("K137").Value = 4892.625
$ws.Range("L137").Value = 8863.875
$ws.Range("M137").Value = -2342.625
$ws.Range("N137").Value = -13963.875
# Row 138
$ws.Range("H138").Value = 3347
$ws.Range("J138").Value = 3356.6
$ws.Range("L138").Value = 10069.8
$ws.Range("N138").Value = -20349.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6927.82
$ws.Range("I32").Value = 4799.8604
$ws.Range("K32").Value = 4799.8604
$ws.Range("M32").Value = -4512.8604
# Row 61
$ws.Range("H61").Value = 2884.8
$ws.Range("I61").Value = 2856.5
$ws.Range("J61").Value = 2998
$ws.Range("K61").Value = 2856.5
$ws.Range("L61").Value = 2998
$ws.Range("M61").Value = -2644.5
$ws.Range("N61").Value = -3422
# Row 97
$ws.Range("H97").Value = 241.27272
$ws.Range("I97").Value = 165.4
$ws.Range("K97").Value = 165.4
$ws.Range("M97").Value = 330.6
# Row 132
$ws.Range("H132").Value = 3022.5881
$ws.Range("I132").Value = 2126.7273
$ws.Range("J132").Value = 4665
$ws.Range("K132").Value = 6380.1819
$ws.Range("L132").Value = 13995
$ws.Range("M132").Value = -3850.1819
$ws.Range("N132").Value = -19055
# Row 136
$ws.Range("H136").Value = 2884.8
$ws.Range("I136").Value = 2856.5
$ws.Range("J136").Value = 2998
$ws.Range("K136").Value = 8569.5
$ws.Range("L136").Value = 8994
$ws.Range("M136").Value = -6019.5
$ws.Range("N136").Value = -14094
# Row 139
$ws.Range("H139").Value = 119999
$ws.Range("J139").Value = 119999
$ws.Range("L139").Value = 119999
$ws.Range("N139").Value = -130279

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 166
$ws.Range("I11").Value = 166
$ws.Range("K11").Value = 166
$ws.Range("M11").Value = -26
# Row 33
$ws.Range("H33").Value = 5000
$ws.Range("I33").Value = 5000
$ws.Range("K33").Value = 5000
$ws.Range("M33").Value = -4664
# Row 64
$ws.Range("H64").Value = 1241.25
$ws.Range("J64").Value = 1319.6666
$ws.Range("L64").Value = 1319.6666
$ws.Range("N64").Value = -1769.6666
# Row 67
$ws.Range("H67").Value = 1241.25
$ws.Range("J67").Value = 1319.6666
$ws.Range("L67").Value = 1319.6666
$ws.Range("N67").Value = -2879.6666
# Row 86
$ws.Range("H86").Value = 2519.4443
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 2519.4443
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
# Row 94
$ws.Range("H94").Value = 1421.375
$ws.Range("I94").Value = 784.1667
$ws.Range("K94").Value = 784.1667
$ws.Range("M94").Value = -333.1667
# Row 105
$ws.Range("H105").Value = 5556.625
$ws.Range("I105").Value = 4857.4
$ws.Range("K105").Value = 4857.4
$ws.Range("M105").Value = -3110.4
# Row 134
$ws.Range("H134").Value = 2507.4
$ws.Range("J134").Value = 2960
$ws.Range("L134").Value = 8880
$ws.Range("N134").Value = -13950
# Row 135
$ws.Range("H135").Value = 81249
$ws.Range("J135").Value = 81249
$ws.Range("L135").Value = 81249
$ws.Range("N135").Value = -91389

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 575.25
$ws.Range("I7").Value = 575.25
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 575.25
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -462.25
# Row 86
$ws.Range("H86").Value = 19087.637
$ws.Range("I86").Value = 9099.091
$ws.Range("J86").Value = 29076.182
$ws.Range("K86").Value = 9099.091
$ws.Range("L86").Value = 29076.182
$ws.Range("M86").Value = -7976.091
$ws.Range("N86").Value = -31322.182
# Row 89
$ws.Range("H89").Value = 19087.637
$ws.Range("I89").Value = 9099.091
$ws.Range("J89").Value = 29076.182
$ws.Range("K89").Value = 45495.455
$ws.Range("L89").Value = 145380.91
$ws.Range("M89").Value = -39879.455
$ws.Range("N89").Value = -156612.91
# Row 105
$ws.Range("H105").Value = 1672.2858
$ws.Range("J105").Value = 1208.5
$ws.Range("L105").Value = 1208.5
$ws.Range("N105").Value = -4702.5
# Row 107
$ws.Range("H107").Value = 1789.4445
$ws.Range("I107").Value = 853
$ws.Range("J107").Value = 2385.3635
$ws.Range("K107").Value = 853
$ws.Range("L107").Value = 2385.3635
$ws.Range("M107").Value = 1067
$ws.Range("N107").Value = -6225.363499999999
# Row 134
$ws.Range("H134").Value = 3273.2
$ws.Range("I134").Value = 3273.2
$ws.Range("K134").Value = 9819.599999999999
$ws.Range("M134").Value = -7284.599999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1333
$ws.Range("I68").Value = 1333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3999
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -3188
# Row 71
$ws.Range("H71").Value = 1333
$ws.Range("I71").Value = 1333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 11997
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -7941
# Row 107
$ws.Range("H107").Value = 2217.6
$ws.Range("J107").Value = 272
$ws.Range("L107").Value = 816
$ws.Range("N107").Value = -4656
# Row 113
$ws.Range("H113").Value = 2123.4
$ws.Range("J113").Value = 2104.3076
$ws.Range("L113").Value = 6312.9228
$ws.Range("N113").Value = -10652.9228

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 668.1818
$ws.Range("I97").Value = 705
$ws.Range("K97").Value = 705
$ws.Range("M97").Value = -209
# Row 136
$ws.Range("H136").Value = 55896.855
$ws.Range("J136").Value = 55896.855
$ws.Range("L136").Value = 167690.565
$ws.Range("N136").Value = -172790.565

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 5249.875
$ws.Range("J22").Value = 5249.875
$ws.Range("L22").Value = 5249.875
$ws.Range("N22").Value = -5839.875
# Row 27
$ws.Range("H27").Value = 5249.875
$ws.Range("J27").Value = 5249.875
$ws.Range("L27").Value = 5249.875
$ws.Range("N27").Value = -5463.875
# Row 55
$ws.Range("H55").Value = 550.5
$ws.Range("I55").Value = 550.5
$ws.Range("K55").Value = 550.5
$ws.Range("M55").Value = -377.5
# Row 93
$ws.Range("H93").Value = 3499.5
$ws.Range("I93").Value = 3499.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3499.5
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -2251.5

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 180
$ws.Range("I113").Value = 142
$ws.Range("J113").Value = 199
$ws.Range("K113").Value = 426
$ws.Range("L113").Value = 597
$ws.Range("M113").Value = 1744
$ws.Range("N113").Value = -4937

Write-Host "Applied Ravana_Profits market data updates"